{"js": "// Replace the 25 division-problem strings in the document's table cells,\n// per the commit diff (each old value is unique and occurs exactly once).\nconst replacements = [[\"804\u00f78=\", \"152\u00f74=\"], [\"684\u00f76=\", \"153\u00f73=\"], [\"655\u00f73=\", \"945\u00f78=\"], [\"783\u00f73=\", \"416\u00f79=\"], [\"295\u00f78=\", \"857\u00f76=\"], [\"433\u00f78=\", \"566\u00f76=\"], [\"629\u00f79=\", \"409\u00f76=\"], [\"194\u00f79=\", \"732\u00f78=\"], [\"251\u00f77=\", \"503\u00f72=\"], [\"705\u00f77=\", \"475\u00f73=\"], [\"453\u00f76=\", \"156\u00f75=\"], [\"971\u00f73=\", \"838\u00f73=\"], [\"114\u00f74=\", \"911\u00f78=\"], [\"373\u00f74=\", \"916\u00f72=\"], [\"839\u00f76=\", \"953\u00f72=\"], [\"468\u00f74=\", \"749\u00f72=\"], [\"491\u00f72=\", \"910\u00f78=\"], [\"249\u00f73=\", \"330\u00f72=\"], [\"312\u00f78=\", \"862\u00f79=\"], [\"316\u00f72=\", \"507\u00f78=\"], [\"239\u00f74=\", \"297\u00f73=\"], [\"420\u00f74=\", \"746\u00f74=\"], [\"707\u00f73=\", \"631\u00f78=\"], [\"641\u00f75=\", \"395\u00f76=\"], [\"914\u00f76=\", \"479\u00f73=\"]];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Each (old, new) pair corresponds to one of the 25 division-problem cells\n# updated by the commit; every old string is unique in the document.\n$pairs = @(\n  @{old=\"804\u00f78=\"; new=\"152\u00f74=\"},\n  @{old=\"684\u00f76=\"; new=\"153\u00f73=\"},\n  @{old=\"655\u00f73=\"; new=\"945\u00f78=\"},\n  @{old=\"783\u00f73=\"; new=\"416\u00f79=\"},\n  @{old=\"295\u00f78=\"; new=\"857\u00f76=\"},\n  @{old=\"433\u00f78=\"; new=\"566\u00f76=\"},\n  @{old=\"629\u00f79=\"; new=\"409\u00f76=\"},\n  @{old=\"194\u00f79=\"; new=\"732\u00f78=\"},\n  @{old=\"251\u00f77=\"; new=\"503\u00f72=\"},\n  @{old=\"705\u00f77=\"; new=\"475\u00f73=\"},\n  @{old=\"453\u00f76=\"; new=\"156\u00f75=\"},\n  @{old=\"971\u00f73=\"; new=\"838\u00f73=\"},\n  @{old=\"114\u00f74=\"; new=\"911\u00f78=\"},\n  @{old=\"373\u00f74=\"; new=\"916\u00f72=\"},\n  @{old=\"839\u00f76=\"; new=\"953\u00f72=\"},\n  @{old=\"468\u00f74=\"; new=\"749\u00f72=\"},\n  @{old=\"491\u00f72=\"; new=\"910\u00f78=\"},\n  @{old=\"249\u00f73=\"; new=\"330\u00f72=\"},\n  @{old=\"312\u00f78=\"; new=\"862\u00f79=\"},\n  @{old=\"316\u00f72=\"; new=\"507\u00f78=\"},\n  @{old=\"239\u00f74=\"; new=\"297\u00f73=\"},\n  @{old=\"420\u00f74=\"; new=\"746\u00f74=\"},\n  @{old=\"707\u00f73=\"; new=\"631\u00f78=\"},\n  @{old=\"641\u00f75=\"; new=\"395\u00f76=\"},\n  @{old=\"914\u00f76=\"; new=\"479\u00f73=\"}\n)\n\nforeach ($p in $pairs) {\n  $range = $d.Content\n  $found = $range.Find.Execute($p.old, $false, $false, $false, $false, $false, $true, 1, $false, $p.new, 2)\n  if (-not $found) {\n    throw \"Find/Replace failed for: \" + $p.old\n  }\n}\n"}
